$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.426.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.89%  '
$ws.Range("D3").Value = "'2.649.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'522.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = "'144.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").Value = "'0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.24%  '
$ws.Range("D9").Value = "'6.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.42%  '
$ws.Range("E10").Value = '  -3.21%  '
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = "'3.117.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.58%  '
$ws.Range("D14").Value = "'58.422.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.83%  '
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = "'2.652.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -9.86%  '
$ws.Range("D18").Value = "'339.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.04%  '
$ws.Range("E19").Value = '  -2.87%  '
$ws.Range("E20").Value = '  -1.23%  '
$ws.Range("D21").Value = "'6.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = "'64.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("E24").Value = '  +0.88%  '
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("E27").Value = '  -2.46%  '
$ws.Range("E28").Value = '  -3.34%  '
$ws.Range("D29").Value = "'6.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.97%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").Value = "'152.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.92%  '
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("D34").Value = "'4.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.76%  '
$ws.Range("E35").Value = '  -5.22%  '
$ws.Range("E36").Value = '  -4.50%  '
$ws.Range("D37").Value = "'0.868"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").Value = "'36.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("E39").Value = '  -4.73%  '
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").Value = "'0.610"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").Value = "'274.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.41%  '
$ws.Range("E44").Value = '  -2.06%  '
$ws.Range("D45").Value = "'19.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.51%  '
$ws.Range("D46").Value = "'0.0536"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").Value = "'10.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.48%  '
$ws.Range("D48").Value = "'2.045.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.75%  '
$ws.Range("D49").Value = "'4.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.25%  '
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("D51").Value = "'18.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.64%  '
